# Apply quarterly financials update to the MGP sheet.
# Two new quarterly columns are inserted before column D, which pushes the
# existing quarters from D:K over to F:M, and the two freshest quarters of
# data are populated into the new D:E columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new (blank) columns before column D; existing D:K data shifts to F:M
$ws.Columns("D:E").Insert()

# Copy the number formats/styles from the (now shifted) F:G columns into the new D:E columns
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D (latest quarter) and column E (prior quarter) with data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 284000
$ws.Range("E8").Value = 282200
$ws.Range("D9:E9").Value = "NA"
$ws.Range("D10:E10").Value = "NA"
$ws.Range("D12:E12").Value = "NA"
$ws.Range("D13:E13").Value = 0
$ws.Range("D14").Value = 1800
$ws.Range("E14").Value = 4400
$ws.Range("D15").Value = 70200
$ws.Range("E15").Value = 66700
$ws.Range("D17").Value = 153300
$ws.Range("E17").Value = 147400
$ws.Range("D18").Value = 130700
$ws.Range("E18").Value = 134800
$ws.Range("D20:E20").Value = -800
$ws.Range("D21").Value = 200100
$ws.Range("E21").Value = 200700
$ws.Range("D22").Value = 58300
$ws.Range("E22").Value = 58700
$ws.Range("D23").Value = 71600
$ws.Range("E23").Value = 75200
$ws.Range("D24").Value = 3100
$ws.Range("E24").Value = 5300
$ws.Range("D25:E25").Value = 0
$ws.Range("D26").Value = 68600
$ws.Range("E26").Value = 69900
$ws.Range("D27").Value = 18600
$ws.Range("E27").Value = 19500
$ws.Range("D28:E28").Value = 0
$ws.Range("D29:E29").Value = 0
$ws.Range("D30:E30").Value = 0
$ws.Range("D31:E31").Value = 0
$ws.Range("D32:E32").Value = 800
$ws.Range("D33").Value = 18600
$ws.Range("E33").Value = 19500
$ws.Range("D34:E34").Value = 0
$ws.Range("D35").Value = 18600
$ws.Range("E35").Value = 19500
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 59800
$ws.Range("E41").Value = 49500
$ws.Range("D42:E42").Value = 0
$ws.Range("D43").Value = 15000
$ws.Range("E43").Value = 12400
$ws.Range("D44:E44").Value = 0
$ws.Range("D45").Value = 37800
$ws.Range("E45").Value = 56400
$ws.Range("D46:E46").Value = 0
$ws.Range("D47:E47").Value = 0
$ws.Range("D48").Value = 10526500
$ws.Range("E48").Value = 10592400
$ws.Range("D49").Value = 312100
$ws.Range("E49").Value = 313400
$ws.Range("D50:E50").Value = 0
$ws.Range("D51:E51").Value = 0
$ws.Range("D52:E52").Value = 0
$ws.Range("D53:E53").Value = 0
$ws.Range("D54").Value = 10951300
$ws.Range("E54").Value = 11024200
$ws.Range("D57").Value = 49600
$ws.Range("E57").Value = 39600
$ws.Range("D58:E58").Value = 0
$ws.Range("D59").Value = 145500
$ws.Range("E59").Value = 149200
$ws.Range("D60:E60").Value = 0
$ws.Range("D61").Value = 4666900
$ws.Range("E61").Value = 4684700
$ws.Range("D62").Value = 33600
$ws.Range("E62").Value = 31400
$ws.Range("D63:E63").Value = 0
$ws.Range("D64:E64").Value = 0
$ws.Range("D65:E65").Value = 0
$ws.Range("D66").Value = 9385300
$ws.Range("E66").Value = 9439800
$ws.Range("D68:E68").Value = 0
$ws.Range("D69:E69").Value = 0
$ws.Range("D70:E70").Value = 0
$ws.Range("D71:E71").Value = 0
$ws.Range("D72").Value = -150900
$ws.Range("E72").Value = -137800
$ws.Range("D73:E73").Value = 0
$ws.Range("D74:E74").Value = 0
$ws.Range("D75:E75").Value = 0
$ws.Range("D76").Value = 1566000
$ws.Range("E76").Value = 1584400
$ws.Range("D77:E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 18600
$ws.Range("E81").Value = 19500
$ws.Range("D83").Value = 70200
$ws.Range("E83").Value = 66700
$ws.Range("D84:E84").Value = 0
$ws.Range("D85:E85").Value = 0
$ws.Range("D86:E86").Value = 0
$ws.Range("D87:E87").Value = 0
$ws.Range("D88:E88").Value = 0
$ws.Range("D89").Value = 147100
$ws.Range("E89").Value = 148800
$ws.Range("D91").Value = -800
$ws.Range("E91").Value = -600
$ws.Range("D92:E92").Value = 0
$ws.Range("D93:E93").Value = 0
$ws.Range("D94").Value = -800
$ws.Range("E94").Value = -1035100
$ws.Range("D96").Value = -116400
$ws.Range("E96").Value = -114400
$ws.Range("D97:E97").Value = 0
$ws.Range("D98:E98").Value = 0
$ws.Range("D99:E99").Value = 0
$ws.Range("D100").Value = -136000
$ws.Range("E100").Value = 646000
$ws.Range("D101:E101").Value = 0
$ws.Range("D102").Value = 10300
$ws.Range("E102").Value = -240400

# Row 91 (Capital Expenditures, Cash Flow section): quarters that shifted into
# H91:J91 were restated from numeric 0/-500 figures to "NA" in this update
$ws.Range("H91:J91").Value = "NA"
